$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.346.07"
$ws.Range("E2").Value = "  +5.13%  "
$ws.Range("D3").Value = "3.171.93"
$ws.Range("E3").Value = "  +2.92%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'398.32"
$ws.Range("E5").Value = "  +2.50%  "
$ws.Range("D6").Value = "'109.20"
$ws.Range("E6").Value = "  +5.41%  "
$ws.Range("D7").Value = "'0.547"
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "'0.616"
$ws.Range("E9").Value = "  +4.81%  "
$ws.Range("D10").Value = "'38.95"
$ws.Range("E10").Value = "  +5.07%  "
$ws.Range("D12").Value = "'0.0879"
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("D13").Value = "3.663.32"
$ws.Range("E13").Value = "  +2.70%  "
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("D15").Value = "'8.03"
$ws.Range("E15").Value = "  +2.77%  "
$ws.Range("E16").Value = "  +8.19%  "
$ws.Range("D17").Value = "3.174.35"
$ws.Range("D18").Value = "'10.47"
$ws.Range("E18").Value = "  -2.59%  "
$ws.Range("D19").Value = "54.243.53"
$ws.Range("E19").Value = "  +4.74%  "
$ws.Range("E20").Value = "  +3.49%  "
$ws.Range("D21").Value = "'12.84"
$ws.Range("E21").Value = "  +3.00%  "
$ws.Range("D22").Value = "0.0₃0982"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").Value = "'71.12"
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("D24").Value = "'271.94"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("E25").Value = "  +2.61%  "
$ws.Range("D26").Value = "'7.99"
$ws.Range("E26").Value = "  -2.87%  "
$ws.Range("D27").Value = "'27.65"
$ws.Range("E27").Value = "  +2.43%  "
$ws.Range("D28").Value = "'7.37"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").Value = "'0.170"
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "'0.112"
$ws.Range("E31").Value = "  +4.07%  "
$ws.Range("D32").Value = "'10.98"
$ws.Range("E32").Value = "  +6.65%  "
$ws.Range("D33").Value = "'0.0501"
$ws.Range("E33").Value = "  +11.42%  "
$ws.Range("D34").Value = "'36.89"
$ws.Range("E34").Value = "  +5.44%  "
$ws.Range("D35").Value = "'2.09"
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("D36").Value = "'50.53"
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("E37").Value = "  +8.64%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("E39").Value = "  +10.61%  "
$ws.Range("D40").Value = "'4.11"
$ws.Range("E40").Value = "  +9.76%  "
$ws.Range("D41").Value = "'0.291"
$ws.Range("E41").Value = "  -1.50%  "
$ws.Range("E42").Value = "  +2.04%  "
$ws.Range("E43").Value = "  +1.00%  "
$ws.Range("E44").Value = "  +3.31%  "
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("D46").Value = "'22.23"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").Value = "'2.07"
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("D49").Value = "2.086.30"
$ws.Range("E49").Value = "  +2.22%  "
$ws.Range("E50").Value = "  +6.85%  "
$ws.Range("D51").Value = "'5.73"
$ws.Range("E51").Value = "  +6.05%  "
